# Update the practice-problem table: "three-digit number × one-digit number".
#
# The worksheet table has 20 rows x 5 columns; only rows 1, 5, 10, 15 and 20
# hold problems (the rows between them are blank answer rows). Every problem
# cell's text is replaced with a new problem, addressed directly by
# (row, column) so there is no ambiguity from repeated/duplicate values.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newValues = @{
  1  = @("571×9=", "216×5=", "863×8=", "431×4=", "838×9=")
  5  = @("260×3=", "822×7=", "523×3=", "180×8=", "827×8=")
  10 = @("123×8=", "414×5=", "601×5=", "642×6=", "874×5=")
  15 = @("282×5=", "349×9=", "975×6=", "224×5=", "213×2=")
  20 = @("834×3=", "877×7=", "496×4=", "721×6=", "891×9=")
}

foreach ($rowIndex in $newValues.Keys) {
  $rowValues = $newValues[$rowIndex]
  for ($col = 1; $col -le 5; $col++) {
    $cell = $t.Cell($rowIndex, $col)
    $cell.Range.Text = $rowValues[$col - 1]
  }
}
